$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3 (shifts existing rows 3..21 down to 4..22),
# pushing the whole table down by one and making room for the new IPO entry.
$ws.Rows("3:3").Insert()

# Fill the newly inserted row 3 with the "신영스팩10호" data.
$ws.Range("A3").Value = "신영스팩10호"
$ws.Range("B3").Value = "2024.01.22~01.23"
$ws.Range("C3").Value = "2,000~2,000"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = 9150
$ws.Range("F3").Value = "신영증권"

# The old last row (originally row 21, "NH스팩30호") is now row 22 after the
# insert above; remove it so the table keeps its original 20-data-row size.
$ws.Rows("22:22").Delete()
